# Applies two small edits described in the commit:
#  1) Slide 1 "TextBox 6": rewrite the Motivation paragraphs -
#     drop the "Motivation:" heading line, tweak the Blake & O'Brien
#     sentence, add a blank line and a new closing sentence.
#  2) Slide 3 "TextBox 7": nudge the textbox 1297 EMU to the right.

$p = $ppt.ActivePresentation

function Get-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Name -eq $name) {
            return $shp
        }
    }
    return $null
}

# ---------------------------------------------------------------
# 1) Slide 1 - "TextBox 6"
# ---------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$tb6 = Get-ShapeByName $slide1 "TextBox 6"
$tr = $tb6.TextFrame.TextRange

$cite = "Blake and O" + [char]0x2019 + "Brien 2016"
$rest = " discovered curtains: a latitudinally narrow and stationary form of precipitation with mysterious physical shape, statistical properties, and origin. "
$closing = "We investigated their statistical properties and began to unravel their mysterious origin."

$tr.Text = $cite + $rest + "`r" + "`r" + $closing

# The whole range inherited size/bold/italic from the old first run
# ("Motivation:", bold, 24pt) - restore the common baseline first ...
$tr.Font.Size = 24
$tr.Font.Bold = $false
$tr.Font.Italic = $false

# ... then re-apply the italic citation.
$tr.Characters(1, $cite.Length).Font.Italic = $true

# ---------------------------------------------------------------
# 2) Slide 3 - "TextBox 7"
# ---------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$tb7 = Get-ShapeByName $slide3 "TextBox 7"

# 259714 EMU / 12700 EMU-per-point, nudged slightly within the
# float32 rounding bucket so the stored EMU lands exactly on 259714.
$tb7.Left = 20.44995
